$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (values such as "27.703.47" or "0.07370" must not be reinterpreted as
# numbers by Excel, which would drop formatting like trailing zeros or
# the multi-dot "thousands" separators used in this sheet).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.703.47"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.847.75"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -2.05%  "
$ws.Range("D5").Value = "320.34"
$ws.Range("E5").Value = "  -1.10%  "
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("D7").Value = "0.4308"
$ws.Range("E7").Value = "  -2.52%  "
$ws.Range("D8").Value = "0.3739"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "0.07370"
$ws.Range("E9").Value = "  -1.24%  "
$ws.Range("D10").Value = "0.8778"
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("D11").Value = "21.65"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "1.854.84"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").Value = "6.725"
$ws.Range("E13").Value = "  -0.51%  "
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "0.07140"
$ws.Range("E15").Value = "  -1.12%  "
$ws.Range("D16").Value = "88.38"
$ws.Range("E16").Value = "  +4.70%  "
$ws.Range("E17").Value = "  -2.07%  "
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").Value = "15.46"
$ws.Range("E20").Value = "  -0.57%  "
$ws.Range("D21").Value = "27.724.52"
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("D22").Value = "5.236"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").Value = "2.089.05"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "2.014"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "155.93"
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  +7.67%  "
$ws.Range("D29").Value = "5.408"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "118.95"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "0.7774"
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "4.567"
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "2.925"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "1.137"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("D38").Value = "0.05353"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "7.272"
$ws.Range("E40").Value = "  +5.98%  "
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").Value = "0.5155"
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("D43").Value = "0.1686"
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "8.837"
$ws.Range("E44").Value = "  +1.92%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "109.50"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "0.4746"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "0.06490"
$ws.Range("D49").Value = "1.695"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "1.014"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "1.851"
$ws.Range("E51").Value = "  -2.93%  "
